# Update cryptocurrency price (D) and volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.000" or
# "24.373.55" are not reinterpreted as numbers and lose their exact formatting
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "24.435.41"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.669.88"
$ws.Range("E3").Value = "  +1.94%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "312.55"
$ws.Range("E5").Value = "  +1.70%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "0.3947"
$ws.Range("E7").Value = "  +1.96%  "
$ws.Range("D8").Value = "0.3923"
$ws.Range("E8").Value = "  +2.55%  "
$ws.Range("D9").Value = "52.12"
$ws.Range("E9").Value = "  +5.74%  "
$ws.Range("D10").Value = "1.390"
$ws.Range("E10").Value = "  +3.62%  "
$ws.Range("D11").Value = "1.003"
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("D12").Value = "0.08572"
$ws.Range("E12").Value = "  +1.00%  "
$ws.Range("D13").Value = "24.51"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "7.286"
$ws.Range("E14").Value = "  +3.36%  "
$ws.Range("D15").Value = "8.019"
$ws.Range("E15").Value = "  +7.98%  "
$ws.Range("D16").Value = "0.00001334"
$ws.Range("E16").Value = "  +4.61%  "
$ws.Range("D17").Value = "1.662.67"
$ws.Range("E17").Value = "  +1.52%  "
$ws.Range("D18").Value = "95.03"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "0.07038"
$ws.Range("E19").Value = "  +2.39%  "
$ws.Range("D20").Value = "20.60"
$ws.Range("E20").Value = "  -0.59%  "
$ws.Range("D21").Value = "6.992"
$ws.Range("E21").Value = "  +1.83%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("D23").Value = "13.75"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").Value = "24.418.92"
$ws.Range("E24").Value = "  +1.41%  "
$ws.Range("D25").Value = "2.504"
$ws.Range("E25").Value = "  +7.61%  "
$ws.Range("D26").Value = "3.100"
$ws.Range("E26").Value = "  +16.03%  "
$ws.Range("D27").Value = "22.55"
$ws.Range("E27").Value = "  +1.33%  "
$ws.Range("D28").Value = "157.13"
$ws.Range("E28").Value = "  +0.21%  "
$ws.Range("D29").Value = "142.87"
$ws.Range("E29").Value = "  +2.30%  "
$ws.Range("D30").Value = "5.348"
$ws.Range("E30").Value = "  +0.12%  "
$ws.Range("D31").Value = "7.948"
$ws.Range("E31").Value = "  -8.28%  "
$ws.Range("D32").Value = "2.555"
$ws.Range("E32").Value = "  +5.97%  "
$ws.Range("D33").Value = "1.845.80"
$ws.Range("E33").Value = "  +1.54%  "
$ws.Range("D34").Value = "1.064"
$ws.Range("E34").Value = "  +12.81%  "
$ws.Range("D35").Value = "0.03130"
$ws.Range("E35").Value = "  +9.14%  "
$ws.Range("D36").Value = "0.08233"
$ws.Range("E36").Value = "  +3.32%  "
$ws.Range("D37").Value = "6.895"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").Value = "11.26"
$ws.Range("E38").Value = "  +14.58%  "
$ws.Range("D39").Value = "0.2767"
$ws.Range("E39").Value = "  +4.11%  "
$ws.Range("D40").Value = "0.09260"
$ws.Range("E40").Value = "  +1.27%  "
$ws.Range("D41").Value = "0.7702"
$ws.Range("E41").Value = "  +2.77%  "
$ws.Range("D42").Value = "13.78"
$ws.Range("E42").Value = "  +6.45%  "
$ws.Range("D43").Value = "1.446"
$ws.Range("E43").Value = "  +0.46%  "
$ws.Range("D44").Value = "16.62"
$ws.Range("E44").Value = "  +4.48%  "
$ws.Range("D45").Value = "0.7096"
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("D46").Value = "2.543"
$ws.Range("E46").Value = "  +4.10%  "
$ws.Range("E47").Value = "  +1.33%  "
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "0.08429"
$ws.Range("E49").Value = "  +1.11%  "
$ws.Range("D50").Value = "136.58"
$ws.Range("E50").Value = "  +3.42%  "
$ws.Range("D51").Value = "1.265"
$ws.Range("E51").Value = "  +0.97%  "
